$wb = $excel.ActiveWorkbook

# --- Roadmap becomes the active/selected sheet tab (was "Models") ---
$ws = $wb.Worksheets.Item("Roadmap")
[void]$ws.Activate()

# Restore/park the window near the top of the sheet and leave the
# selection on C14, matching the saved view state of the edited file.
$win = $excel.ActiveWindow
$win.Left = 880
$win.Top = 0
$win.Width = 25600
$win.Height = 13200

# --- Fill in the newly-added "Ensemble Model" row (row 13) of the
# Week 1 table. Cells are populated in this order (Results/"presentation"
# columns first, then the Notes column, then the Model Name/Achieved
# columns) to mirror how the new text entries were authored. ---
$ws.Range("E13").Value = "presentation / insights / slides"
$ws.Range("F13").Value = "presentation / slides"
$ws.Range("D13").Value = "modeling" + [char]10 + "user testing"
$ws.Range("B13").Value = "Word2vec model" + [char]10 + "NMF Model"
$ws.Range("C13").Value = "Ensemble Model"

# B13/D13 hold multi-line text, so wrap them like the other
# model-name / notes cells in this table (e.g. B6, B7, F14).
$ws.Range("B13").WrapText = $true
$ws.Range("D13").WrapText = $true

# Leave the cursor/selection on C14 of the Roadmap sheet.
[void]$ws.Range("C14").Select()
